$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The comma-separated bounding box coordinates (column I) and the
# confidence value (column J) look numeric to Excel, so the cells are
# formatted as Text first to keep them stored as strings (matching the
# original inline-string cell values) rather than being reinterpreted as
# numbers.
$ws.Range("I16").NumberFormat = "@"
$ws.Range("I17").NumberFormat = "@"
$ws.Range("J16").NumberFormat = "@"

# Row 16 updates
$ws.Range("D16").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("I16").Value = "642,530,686,576"
$ws.Range("J16").Value = "0.75"

# Row 17 updates
$ws.Range("D17").Value = "image_20250807111344_ppp0.jpg"
$ws.Range("I17").Value = "794,481,831,526"
